# Adds a "medicaid" boolean column (E) to the ASC_sub_2 sheet, mirroring
# the header styling of the existing A1:D1 header row, and populating
# E2:E67 with the per-patient medicaid flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (E1) ---
$ws.Range("E1").Value = "medicaid"

# Copy the header formatting (bold font, border, centered/top alignment)
# from the existing D1 header cell onto the new E1 header cell.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data cells (E2:E67) ---
# medicaid flag for each patient row, in row order (row 2 .. row 67)
$medicaid = @(
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,
    0,1,0,0,0,0,0,0,0,0,
    0,0,0,1,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,1,0,
    1,0,0,0,0,0
)

for ($i = 0; $i -lt $medicaid.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = [bool]$medicaid[$i]
}
